$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 4999
$ws.Range("J13").Value = 5331.6665
$ws.Range("L13").Value = 5331.6665
$ws.Range("N13").Value = -5669.6665
$ws.Range("H32").Value = 2250
$ws.Range("I32").Value = 1500
$ws.Range("K32").Value = 1500
$ws.Range("M32").Value = -1174
$ws.Range("H40").Value = 4166
$ws.Range("J40").Value = 4748.5
$ws.Range("L40").Value = 4748.5
$ws.Range("N40").Value = -5098.5
$ws.Range("H53").Value = 129.75
$ws.Range("I53").Value = 57.857143
$ws.Range("J53").Value = 230.4
$ws.Range("K53").Value = 57.857143
$ws.Range("L53").Value = 230.4
$ws.Range("M53").Value = 579.142857
$ws.Range("N53").Value = -1504.4
$ws.Range("H125").Value = 30090.285
$ws.Range("J125").Value = 34822.168
$ws.Range("L125").Value = 313399.512
$ws.Range("N125").Value = -318319.512
$ws.Range("H130").Value = 70000
$ws.Range("J130").Value = 70000
$ws.Range("L130").Value = 70000
$ws.Range("N130").Value = -80040
$ws.Range("H131").Value = 1500
$ws.Range("I131").Value = 1500
$ws.Range("K131").Value = 4500
$ws.Range("M131").Value = 540
$ws.Range("H132").Value = 2071.348
$ws.Range("I132").Value = 1302.4706
$ws.Range("K132").Value = 3907.4118
$ws.Range("M132").Value = -1377.4118
$ws.Range("H135").Value = 814.6667
$ws.Range("I135").Value = 837.08105
$ws.Range("K135").Value = 7533.72945
$ws.Range("M135").Value = -4998.72945
$ws.Range("H137").Value = 2881.0715
$ws.Range("I137").Value = 2576.4285
$ws.Range("J137").Value = 3185.7144
$ws.Range("K137").Value = 7729.2855
$ws.Range("L137").Value = 9557.143199999999
$ws.Range("M137").Value = -5179.2855
$ws.Range("N137").Value = -14657.1432
$ws.Range("H138").Value = 3852.5066
$ws.Range("I138").Value = 1831.7
$ws.Range("J138").Value = 5142.383
$ws.Range("K138").Value = 5495.1
$ws.Range("L138").Value = 15427.149
$ws.Range("M138").Value = -355.1000000000004
$ws.Range("N138").Value = -25707.149

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H12").Value = 2333.3333
$ws.Range("J12").Value = 2000
$ws.Range("L12").Value = 2000
$ws.Range("N12").Value = -2346
$ws.Range("H32").Value = 3673.0908
$ws.Range("I32").Value = 3188.6904
$ws.Range("K32").Value = 3188.6904
$ws.Range("M32").Value = -2901.6904
$ws.Range("H74").Value = 1822.45
$ws.Range("I74").Value = 1860.4736
$ws.Range("K74").Value = 1860.4736
$ws.Range("M74").Value = -986.4736
$ws.Range("H77").Value = 1822.45
$ws.Range("I77").Value = 1860.4736
$ws.Range("K77").Value = 9302.368
$ws.Range("M77").Value = -4934.368

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1607.1296
$ws.Range("I134").Value = 1386.9131
$ws.Range("K134").Value = 4160.7393
$ws.Range("M134").Value = -1625.7393

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 13499.333
$ws.Range("J9").Value = 13499.333
$ws.Range("L9").Value = 13499.333
$ws.Range("N9").Value = -13835.333
$ws.Range("H22").Value = 660
$ws.Range("I22").Value = 443.44446
$ws.Range("J22").Value = 774.64703
$ws.Range("K22").Value = 443.44446
$ws.Range("L22").Value = 774.64703
$ws.Range("M22").Value = -93.44445999999999
$ws.Range("N22").Value = -1474.64703
$ws.Range("H35").Value = 908.875
$ws.Range("I35").Value = 155.6
$ws.Range("J35").Value = 2164.3333
$ws.Range("K35").Value = 155.6
$ws.Range("L35").Value = 2164.3333
$ws.Range("M35").Value = 138.4
$ws.Range("N35").Value = -2752.3333
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("H58").Value = 2262.5925
$ws.Range("I58").Value = 3104.077
$ws.Range("K58").Value = 3104.077
$ws.Range("M58").Value = -2901.077
$ws.Range("H61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()
$ws.Range("H86").Value = 6174.9
$ws.Range("I86").Value = 6247.25
$ws.Range("K86").Value = 6247.25
$ws.Range("M86").Value = -5124.25
$ws.Range("H89").Value = 6174.9
$ws.Range("I89").Value = 6247.25
$ws.Range("K89").Value = 31236.25
$ws.Range("M89").Value = -25620.25
$ws.Range("H134").Value = 1501.5555
$ws.Range("I134").Value = 1411.9183
$ws.Range("K134").Value = 4235.7549
$ws.Range("M134").Value = -1700.7549
$ws.Range("H136").Value = 2262.5925
$ws.Range("I136").Value = 3104.077
$ws.Range("K136").Value = 9312.231
$ws.Range("M136").Value = -6762.231

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 144.33333
$ws.Range("I7").Value = 153
$ws.Range("J7").Value = 118.333336
$ws.Range("K7").Value = 459
$ws.Range("L7").Value = 355.000008
$ws.Range("M7").Value = -347
$ws.Range("N7").Value = -579.000008
$ws.Range("H92").Value = 1033.3334
$ws.Range("I92").Value = 500
$ws.Range("J92").Value = 1300
$ws.Range("K92").Value = 1500
$ws.Range("L92").Value = 3900
$ws.Range("M92").Value = -252
$ws.Range("N92").Value = -6396

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 5889.933
$ws.Range("I122").Value = 10045.917
$ws.Range("J122").Value = 3119.2778
$ws.Range("K122").Value = 30137.751
$ws.Range("L122").Value = 9357.8334
$ws.Range("M122").Value = -27687.751
$ws.Range("N122").Value = -14257.8334

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1050.9565
$ws.Range("I16").Value = 551.2353000000001
$ws.Range("J16").Value = 2466.8333
$ws.Range("K16").Value = 551.2353000000001
$ws.Range("L16").Value = 2466.8333
$ws.Range("M16").Value = -381.2353000000001
$ws.Range("N16").Value = -2806.8333
$ws.Range("H46").Value = 3322.1667
$ws.Range("I46").Value = 900
$ws.Range("J46").Value = 3542.3635
$ws.Range("K46").Value = 900
$ws.Range("L46").Value = 3542.3635
$ws.Range("M46").Value = -712
$ws.Range("N46").Value = -3918.3635
$ws.Range("H61").Value = 3576.617
$ws.Range("I61").Value = 2819.2144
$ws.Range("K61").Value = 2819.2144
$ws.Range("M61").Value = -2617.2144
$ws.Range("H113").Value = 3576.617
$ws.Range("I113").Value = 2819.2144
$ws.Range("K113").Value = 2819.2144
$ws.Range("M113").Value = -649.2143999999998
$ws.Range("H122").Value = 7526.923
$ws.Range("J122").Value = 11385.333
$ws.Range("L122").Value = 34155.999
$ws.Range("N122").Value = -39055.999
$ws.Range("H127").Value = 222000
$ws.Range("J127").Value = 222000
$ws.Range("L127").Value = 222000
$ws.Range("N127").Value = -231920
$ws.Range("H132").Value = 2258.1924
$ws.Range("I132").Value = 2193.7856
$ws.Range("J132").Value = 2333.3333
$ws.Range("K132").Value = 6581.3568
$ws.Range("L132").Value = 6999.999899999999
$ws.Range("M132").Value = -4051.3568
$ws.Range("N132").Value = -12059.9999
$ws.Range("H136").Value = 19454.066
$ws.Range("I136").Value = 1468.8572
$ws.Range("K136").Value = 4406.571599999999
$ws.Range("M136").Value = -1856.571599999999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 644.96295
$ws.Range("I113").Value = 365.72726
$ws.Range("K113").Value = 1097.18178
$ws.Range("M113").Value = 1072.81822
$ws.Range("H122").Value = 5184
$ws.Range("I122").Value = 1979.4
$ws.Range("J122").Value = 6328.5
$ws.Range("K122").Value = 5938.200000000001
$ws.Range("L122").Value = 18985.5
$ws.Range("M122").Value = -3488.200000000001
$ws.Range("N122").Value = -23885.5
$ws.Range("H132").Value = 4768.8184
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("H136").Value = 2162.7
$ws.Range("I136").Value = 1600.25
$ws.Range("J136").Value = 3006.375
$ws.Range("K136").Value = 4800.75
$ws.Range("L136").Value = 9019.125
$ws.Range("M136").Value = -2250.75
$ws.Range("N136").Value = -14119.125
